# Apply updated cryptos data (price/volume/hora) for Fri Feb 3 2023 09:09:39 UTC run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "323.46"
Set-TextValue "E2" "-2.06%"
Set-TextValue "G2" "9"
Set-TextValue "D3" "39.35"
Set-TextValue "E3" "-1.53%"
Set-TextValue "G3" "9"
Set-TextValue "D4" "5.707"
Set-TextValue "E4" "8.34%"
Set-TextValue "G4" "9"
Set-TextValue "D5" "0.08000"
Set-TextValue "E5" "-1.32%"
Set-TextValue "G5" "9"
Set-TextValue "E6" "-0.50%"
Set-TextValue "G6" "9"
Set-TextValue "D7" "8.611"
Set-TextValue "E7" "-0.06%"
Set-TextValue "G7" "9"
Set-TextValue "D8" "1.980"
Set-TextValue "E8" "2.31%"
Set-TextValue "G8" "9"
Set-TextValue "D9" "2.952"
Set-TextValue "E9" "-0.85%"
Set-TextValue "G9" "9"
Set-TextValue "D10" "0.9277"
Set-TextValue "E10" "-0.94%"
Set-TextValue "G10" "9"
Set-TextValue "D11" "0.1244"
Set-TextValue "E11" "-5.94%"
Set-TextValue "G11" "9"
Set-TextValue "D12" "0.1956"
Set-TextValue "E12" "-0.51%"
Set-TextValue "G12" "9"
Set-TextValue "D13" "8.703"
Set-TextValue "E13" "24.96%"
Set-TextValue "G13" "9"
Set-TextValue "D14" "0.09207"
Set-TextValue "E14" "-0.58%"
Set-TextValue "G14" "9"
Set-TextValue "D15" "0.03585"
Set-TextValue "E15" "0.55%"
Set-TextValue "G15" "9"
Set-TextValue "D16" "0.1049"
Set-TextValue "E16" "9.55%"
Set-TextValue "G16" "9"
Set-TextValue "D17" "0.001303"
Set-TextValue "E17" "-2.36%"
Set-TextValue "G17" "9"
Set-TextValue "D18" "0.006144"
Set-TextValue "E18" "0.11%"
Set-TextValue "G18" "9"
Set-TextValue "D19" "3.349"
Set-TextValue "E19" "-0.61%"
Set-TextValue "G19" "9"
Set-TextValue "D20" "0.3530"
Set-TextValue "E20" "0.15%"
Set-TextValue "G20" "9"
Set-TextValue "D21" "0.1371"
Set-TextValue "E21" "3.73%"
Set-TextValue "G21" "9"
Set-TextValue "D22" "0.2413"
Set-TextValue "E22" "-5.78%"
Set-TextValue "G22" "9"
Set-TextValue "D23" "0.04416"
Set-TextValue "E23" "-0.19%"
Set-TextValue "G23" "9"
Set-TextValue "D24" "0.001264"
Set-TextValue "E24" "3.38%"
Set-TextValue "G24" "9"
Set-TextValue "D25" "0.004606"
Set-TextValue "E25" "6.81%"
Set-TextValue "G25" "9"
Set-TextValue "D26" "0.0001151"
Set-TextValue "E26" "-3.32%"
Set-TextValue "G26" "9"
Set-TextValue "G27" "9"
Set-TextValue "G28" "9"
Set-TextValue "G29" "9"
Set-TextValue "G30" "9"
Set-TextValue "G31" "9"
Set-TextValue "G32" "9"
Set-TextValue "G33" "9"
Set-TextValue "G34" "9"
Set-TextValue "G35" "9"
Set-TextValue "G36" "9"
Set-TextValue "G37" "9"
Set-TextValue "G38" "9"
Set-TextValue "D39" "0.02507"
Set-TextValue "E39" "-0.23%"
Set-TextValue "G39" "9"
Set-TextValue "D40" "0.05348"
Set-TextValue "E40" "3.71%"
Set-TextValue "G40" "9"
Set-TextValue "D41" "0.007472"
Set-TextValue "E41" "-2.65%"
Set-TextValue "G41" "9"
Set-TextValue "D42" "0.009612"
Set-TextValue "E42" "3.92%"
Set-TextValue "G42" "9"
Set-TextValue "D43" "0.1405"
Set-TextValue "E43" "-1.63%"
Set-TextValue "G43" "9"
Set-TextValue "D44" "0.002118"
Set-TextValue "E44" "-2.46%"
Set-TextValue "G44" "9"
Set-TextValue "D45" "0.01017"
Set-TextValue "E45" "-0.72%"
Set-TextValue "G45" "9"
Set-TextValue "D46" "0.00006753"
Set-TextValue "E46" "1.49%"
Set-TextValue "G46" "9"
Set-TextValue "E47" "0.05%"
Set-TextValue "G47" "9"
Set-TextValue "E48" "-11.12%"
Set-TextValue "G48" "9"
Set-TextValue "D49" "0.002292"
Set-TextValue "E49" "-7.67%"
Set-TextValue "G49" "9"
Set-TextValue "D50" "0.00002102"
Set-TextValue "E50" "0.05%"
Set-TextValue "G50" "9"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "0.05%"
Set-TextValue "G51" "9"
